# Update "想去人数" (interest count) figures in column F across all four
# sheets of the workbook, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12408
$ws1.Range("F3").Value  = 6958
$ws1.Range("F6").Value  = 439
$ws1.Range("F9").Value  = 14
$ws1.Range("F10").Value = 972
$ws1.Range("F11").Value = 124
$ws1.Range("F12").Value = 326
$ws1.Range("F13").Value = 980
$ws1.Range("F14").Value = 3708
$ws1.Range("F18").Value = 217
$ws1.Range("F21").Value = 256
$ws1.Range("F22").Value = 288
$ws1.Range("F23").Value = 27
$ws1.Range("F24").Value = 93
$ws1.Range("F25").Value = 341
$ws1.Range("F26").Value = 5140
$ws1.Range("F27").Value = 63
$ws1.Range("F28").Value = 1358
$ws1.Range("F29").Value = 275
$ws1.Range("F30").Value = 841
$ws1.Range("F31").Value = 1295

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 3731
$ws2.Range("F13").Value = 9
$ws2.Range("F17").Value = 44

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9195
$ws3.Range("F3").Value = 540
$ws3.Range("F4").Value = 1931

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9195
$ws4.Range("F3").Value  = 540
$ws4.Range("F4").Value  = 1931
$ws4.Range("F5").Value  = 12408
$ws4.Range("F6").Value  = 6959
$ws4.Range("F8").Value  = 3731
$ws4.Range("F11").Value = 439
$ws4.Range("F14").Value = 14
$ws4.Range("F15").Value = 972
$ws4.Range("F16").Value = 124
$ws4.Range("F17").Value = 326
$ws4.Range("F18").Value = 980
$ws4.Range("F19").Value = 3708
$ws4.Range("F22").Value = 217
$ws4.Range("F25").Value = 256
$ws4.Range("F26").Value = 288
$ws4.Range("F27").Value = 27
$ws4.Range("F32").Value = 341
$ws4.Range("F33").Value = 5140
$ws4.Range("F34").Value = 63
$ws4.Range("F35").Value = 1358
$ws4.Range("F38").Value = 275
$ws4.Range("F40").Value = 841
$ws4.Range("F41").Value = 1295
$ws4.Range("F43").Value = 9
$ws4.Range("F48").Value = 44
